$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Hora column (G2:G51) from 18 to 19 for all data rows
$ws.Range("G2:G51").Value = "19"

# Update other changed cells (Price, Coin, Link, Volume columns)
$ws.Range("D2").Value = "242.58"
$ws.Range("D3").Value = "23.07"
$ws.Range("D4").Value = "5.417"
$ws.Range("D5").Value = "0.05883"
$ws.Range("D6").Value = "3.435"
$ws.Range("D7").Value = "6.549"
$ws.Range("D8").Value = "0.8108"
$ws.Range("D9").Value = "0.9396"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1417"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07423"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03326"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03053"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09333"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.857"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001578"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04687"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005926"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "0.005866"
$ws.Range("E20").Value = "19BitKanKANBestin24h"
$ws.Range("D21").Value = "0.004886"
$ws.Range("D23").Value = "3.566"
$ws.Range("D24").Value = "2.120"
$ws.Range("D25").Value = "0.3221"
$ws.Range("D27").Value = "0.0002287"
$ws.Range("D40").Value = "0.03962"
$ws.Range("D41").Value = "0.006185"
$ws.Range("D44").Value = "0.009088"
$ws.Range("D45").Value = "0.00005186"
$ws.Range("D47").Value = "0.6707"
